$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = -0.2132
$ws.Range("F5").Value = -0.2486

$ws.Range("E9").Value = 0.0203
$ws.Range("F9").Value = -0.007
$ws.Range("G9").Value = -0.0063
$ws.Range("H9").Value = -0.1352
$ws.Range("J9").Value = -0.044
$ws.Range("M9").Value = -0.0315

$ws.Range("E12").Value = 0.5054
$ws.Range("F12").Value = -0.7273
$ws.Range("G12").Value = 0.0803
$ws.Range("I12").Value = -0.5601
$ws.Range("J12").Value = -0.0431
$ws.Range("K12").Value = -0.0266

$ws.Range("E13").Value = -0.1205
$ws.Range("F13").Value = -0.1116
$ws.Range("G13").Value = -0.1283
$ws.Range("H13").Value = -0.1173
$ws.Range("I13").Value = -0.0493
$ws.Range("M13").Value = -0.3036

$ws.Range("E15").Value = -1.7426
$ws.Range("F15").Value = -1.7903
$ws.Range("G15").Value = -3.227
$ws.Range("H15").Value = -2.8697
$ws.Range("I15").Value = -2.1116
$ws.Range("J15").Value = -1.9137
$ws.Range("K15").Value = -3.0431
$ws.Range("M15").Value = -0.6396

$ws.Range("E17").Value = -0.2039
$ws.Range("F17").Value = -0.1816

$ws.Range("E19").Value = -1.1747
$ws.Range("F19").Value = 0.0014
$ws.Range("G19").Value = -0.0335
$ws.Range("H19").Value = -0.0353
$ws.Range("I19").Value = -0.0348
$ws.Range("J19").Value = -0.0224
$ws.Range("K19").Value = -1.0898

$ws.Range("E21").Value = 0.7173
$ws.Range("F21").Value = 0.552

$ws.Range("E32").Value = -0.0041
$ws.Range("F32").Value = -0.1463

$ws.Range("E36").Value = 0.0142
$ws.Range("F36").Value = -0.004
$ws.Range("G36").Value = -0.0001
$ws.Range("H36").Value = -0.0048
$ws.Range("K36").Value = 0.0001

$ws.Range("E39").Value = 0.0325
$ws.Range("F39").Value = 0.2237
$ws.Range("G39").Value = 0.2185
$ws.Range("H39").Value = 0.1677
$ws.Range("I39").Value = 0.3438
$ws.Range("J39").Value = 0.1536
$ws.Range("K39").Value = 0.1506

$ws.Range("E40").Value = 0.0415
$ws.Range("F40").Value = 0.0315
$ws.Range("G40").Value = 0.0293
$ws.Range("H40").Value = 0.0269
$ws.Range("I40").Value = -0.0112
$ws.Range("J40").Value = -0.0061
$ws.Range("K40").Value = -0.0044
$ws.Range("M40").Value = -0.0028

$ws.Range("E42").Value = 0.2458
$ws.Range("F42").Value = 0.1578
$ws.Range("G42").Value = 0.0838
$ws.Range("H42").Value = -0.1102
$ws.Range("I42").Value = 0.1743
$ws.Range("J42").Value = -0.0341
$ws.Range("K42").Value = -0.0455
$ws.Range("M42").Value = -0.1247

$ws.Range("E44").Value = -0.0647
$ws.Range("F44").Value = -0.1457

$ws.Range("E46").Value = 0.0368
$ws.Range("F46").Value = -0.0104
$ws.Range("G46").Value = -0.0348
$ws.Range("H46").Value = -0.0023
$ws.Range("I46").Value = 0.0001
$ws.Range("J46").Value = 0.0121

$ws.Range("E48").Value = 0.0658
$ws.Range("F48").Value = 0.1467
